# Add five new "Название файла в скриптах и цвет" rows (28-32) to the
# bottom of the table on the active sheet, mirroring the formatting of
# the last existing data row (27) and referencing five new shared
# strings with the script/file names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    "SCRIPT/P02P01A/us0203.ssb",
    "SCRIPT/P02P01A/us0403.ssb",
    "SCRIPT/P02P01A/us2003.ssb",
    "SCRIPT/P02P01A/us2006.ssb",
    "SCRIPT/P02P01A/us2009.ssb"
)

$firstNewRow = 28
$lastNewRow = $firstNewRow + $newValues.Length - 1

# Copy the formatting (style, borders, wrap, etc.) of the last existing
# row's first column cell onto the new range in one shot.
$ws.Range("A27").Copy() | Out-Null
$ws.Range("A$firstNewRow`:A$lastNewRow").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

# Fill in the values and restore the wrapped-text row height used by
# every other data row in the sheet.
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $r = $firstNewRow + $i
    $ws.Range("A$r").Value = $newValues[$i]
    $ws.Rows($r).RowHeight = 43.2
}

# Match the author's final view/selection state.
$win = $excel.ActiveWindow
$ws.Range("D28").Select() | Out-Null
$win.ScrollRow = 25
$win.ScrollColumn = 1

Write-Host "Added rows $firstNewRow-$lastNewRow with $($newValues.Length) new script names."
